$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date and Count values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-05-18T17:38:26+00:00"

# Count's new value ("1") looks like a plain number, so force it to stay text
# (matching the original shared-string cell type) via the classic leading
# apostrophe, then restore the wrap-text formatting that the text-literal
# entry otherwise drops.
$countCell = $meta.Range("B21")
$countCell.Value = "'1"
$countCell.WrapText = $true

# --- Concepts sheet: collapse the concept table down to a single row ---
$concepts = $wb.Worksheets.Item("Concepts")

# Remove the extra data rows (3-7), leaving only the header row and one data row.
$concepts.Rows("3:7").Delete()

# Replace the remaining data row's contents with the new single concept.
# Column A ("Level") already holds the text "1" and doesn't need to change.
$concepts.Range("B2").Value = "GENO"
$concepts.Range("C2").Value = "Genomics"
$concepts.Range("D2").Value = ""
